$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.445.01"
$ws.Range("E2").Value = "  -0.63%  "
$ws.Range("D3").Value = "2.575.28"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("D5").Value = "'581.85"
$ws.Range("E5").Value = "  -1.32%  "
$ws.Range("D6").Value = "'165.16"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'0.526"
$ws.Range("E8").Value = "  -1.43%  "
$ws.Range("D9").Value = "2.574.00"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("E11").Value = "  +0.35%  "
$ws.Range("E12").Value = "  -0.49%  "
$ws.Range("D13").Value = "'5.16"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "'26.76"
$ws.Range("E14").Value = "  -2.59%  "
$ws.Range("D15").Value = "3.041.47"
$ws.Range("E15").Value = "  -2.43%  "
$ws.Range("E16").Value = "  -1.71%  "
$ws.Range("D17").Value = "66.319.65"
$ws.Range("E17").Value = "  -0.76%  "
$ws.Range("D18").Value = "2.570.72"
$ws.Range("E18").Value = "  -2.21%  "
$ws.Range("D19").Value = "'11.42"
$ws.Range("E19").Value = "  -4.15%  "
$ws.Range("D20").Value = "'7.73"
$ws.Range("E20").Value = "  -3.58%  "
$ws.Range("D21").Value = "'351.07"
$ws.Range("E21").Value = "  -1.71%  "
$ws.Range("D22").Value = "'4.23"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").Value = "'4.59"
$ws.Range("E23").Value = "  -2.13%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("B25").Value = "SuiNetwork"
$ws.Range("C25").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D25").Value = "'1.88"
$ws.Range("E25").Value = "  -3.69%  "
$ws.Range("B26").Value = "Aptos"
$ws.Range("C26").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").Value = "'10.04"
$ws.Range("E26").Value = "  -8.04%  "
$ws.Range("D27").Value = "'68.96"
$ws.Range("E27").Value = "  -2.27%  "
$ws.Range("D28").Value = "2.707.92"
$ws.Range("E28").Value = "  -1.73%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "0.0₃0986"
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").Value = "'535.23"
$ws.Range("E31").Value = "  -2.45%  "
$ws.Range("D32").Value = "'8.04"
$ws.Range("E32").Value = "  +1.54%  "
$ws.Range("D33").Value = "'1.33"
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("E35").Value = "  -0.69%  "
$ws.Range("D36").Value = "'0.999"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("E37").Value = "  -3.58%  "
$ws.Range("D38").Value = "'156.91"
$ws.Range("E38").Value = "  -0.25%  "
$ws.Range("D39").Value = "'18.72"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("E40").Value = "  -1.52%  "
$ws.Range("E41").Value = "  +1.93%  "
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("D43").Value = "'5.10"
$ws.Range("E43").Value = "  -1.97%  "
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'2.39"
$ws.Range("E45").Value = "  -2.85%  "
$ws.Range("D46").Value = "0.0₆0287"
$ws.Range("E46").Value = "  -2.54%  "
$ws.Range("D47").Value = "'148.82"
$ws.Range("E47").Value = "  -1.93%  "
$ws.Range("D48").Value = "'0.566"
$ws.Range("E48").Value = "  -3.39%  "
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("D50").Value = "'1.69"
$ws.Range("E50").Value = "  -1.66%  "
$ws.Range("D51").Value = "'0.0760"
$ws.Range("E51").Value = "  -1.74%  "
